# Apply cryptos list update (Tue Dec  5 04:31:12 UTC 2023, GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure every updated cell stays plain text (column B/C/D/E are all
# text-formatted in this sheet, e.g. "1.00" / "2.80" / "41.773.83" must
# not be auto-coerced into numbers by the COM Range.Value setter).

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "41.768.85"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +2.60%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.228.18"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +0.40%  "

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.16%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "231.37"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  +1.41%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.623"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  -0.80%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "60.58"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  -6.56%  "

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  +0.05%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.405"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  -0.51%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "58.28"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -1.55%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0901"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +2.34%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.103"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  -0.42%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "2.562.10"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +0.59%  "

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -2.78%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "22.73"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  +1.61%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.801"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -2.78%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "5.62"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -0.24%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.241.25"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  +0.89%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "41.708.54"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +2.48%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0₃0907"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +0.36%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "72.42"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -2.02%  "

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -0.90%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "247.94"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -1.96%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -0.13%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -0.65%  "

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +1.34%  "

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  +0.18%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "169.31"
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -2.17%  "

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -2.52%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "19.91"
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -2.02%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.40"
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -2.78%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.62"
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -7.94%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.122"
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  -1.55%  "

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +4.33%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.69"
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +0.49%  "

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +3.15%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.56"
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  -9.17%  "

$c = $ws.Range("B38")
$c.NumberFormat = "@"
$c.Value = "RenderToken"
$c = $ws.Range("C38")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.61"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -5.67%  "

$c = $ws.Range("B39")
$c.NumberFormat = "@"
$c.Value = "LidoDAOToken"
$c = $ws.Range("C39")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.38"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -3.59%  "

$c = $ws.Range("B40")
$c.NumberFormat = "@"
$c.Value = "BinanceUSD"
$c = $ws.Range("C40")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  +0.26%  "

$c = $ws.Range("B41")
$c.NumberFormat = "@"
$c.Value = "TerraClassic"
$c = $ws.Range("C41")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.000235"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  +6.80%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0239"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  +1.61%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "8.56"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -1.97%  "

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  -1.59%  "

$c = $ws.Range("B45")
$c.NumberFormat = "@"
$c.Value = "Aave"
$c = $ws.Range("C45")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "98.48"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -3.28%  "

$c = $ws.Range("B46")
$c.NumberFormat = "@"
$c.Value = "FTXToken"
$c = $ws.Range("C46")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.46"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  -8.82%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0957"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  +1.89%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.471.20"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  -2.75%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "16.56"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -5.29%  "

$c = $ws.Range("B50")
$c.NumberFormat = "@"
$c.Value = "HuobiToken"
$c = $ws.Range("C50")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.80"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -2.00%  "

$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = "NEARProtocol"
$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.29"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +8.18%  "

"Updated 106 cells"
